$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update arrival count (B2) and local minima count (D2)
$ws.Range("B2").Value = 32
$ws.Range("D2").Value = 8

# Row 5: update arrival probability (B5) and minima probability (D5)
$ws.Range("B5").Value = 0.8
$ws.Range("D5").Value = 0.2
